# edit.ps1 - applies the "Add files via upload" games-sheet update
# (2023 NFL season weeks 7-10 results) to the games workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1. Fill in game results for weeks 7-10 (rows 367-423).
#    Columns: A=season B=week C=date D=team1 E=team2 F=score1 G=score2
# ---------------------------------------------------------------
$games = @()
$games += ,@(367, 2023, 7, 45219, "JAX", "NO", 31, 24)
$games += ,@(368, 2023, 7, 45222, "LV", "CHI", 12, 30)
$games += ,@(369, 2023, 7, 45222, "DET", "BAL", 6, 38)
$games += ,@(370, 2023, 7, 45222, "CLE", "IND", 39, 38)
$games += ,@(371, 2023, 7, 45222, "WAS", "NYG", 7, 14)
$games += ,@(372, 2023, 7, 45222, "BUF", "NE", 25, 29)
$games += ,@(373, 2023, 7, 45222, "ATL", "TB", 16, 13)
$games += ,@(374, 2023, 7, 45222, "ARI", "SEA", 10, 20)
$games += ,@(375, 2023, 7, 45222, "PIT", "LA", 24, 17)
$games += ,@(376, 2023, 7, 45222, "GB", "DEN", 17, 19)
$games += ,@(377, 2023, 7, 45222, "LAC", "KC", 17, 31)
$games += ,@(378, 2023, 7, 45222, "MIA", "PHI", 17, 31)
$games += ,@(379, 2023, 7, 45223, "SF", "MIN", 17, 22)
$games += ,@(380, 2023, 8, 45225, "TB", "BUF", 18, 24)
$games += ,@(381, 2023, 8, 45228, "LA", "DAL", 20, 43)
$games += ,@(382, 2023, 8, 45228, "MIN", "GB", 24, 10)
$games += ,@(383, 2023, 8, 45228, "ATL", "TEN", 23, 28)
$games += ,@(384, 2023, 8, 45228, "NO", "IND", 38, 27)
$games += ,@(385, 2023, 8, 45228, "NE", "MIA", 17, 31)
$games += ,@(386, 2023, 8, 45228, "NYJ", "NYG", 13, 10)
$games += ,@(387, 2023, 8, 45228, "JAX", "PIT", 20, 10)
$games += ,@(388, 2023, 8, 45228, "PHI", "WAS", 38, 31)
$games += ,@(389, 2023, 8, 45228, "HOU", "CAR", 13, 15)
$games += ,@(390, 2023, 8, 45228, "CLE", "SEA", 20, 24)
$games += ,@(391, 2023, 8, 45228, "CIN", "SF", 31, 17)
$games += ,@(392, 2023, 8, 45228, "BAL", "ARI", 31, 24)
$games += ,@(393, 2023, 8, 45228, "KC", "DEN", 9, 24)
$games += ,@(394, 2023, 8, 45228, "CHI", "LAC", 13, 30)
$games += ,@(395, 2023, 8, 45229, "LV", "DET", 14, 26)
$games += ,@(396, 2023, 9, 45232, "TEN ", "PIT", 26, 20)
$games += ,@(397, 2023, 9, 45235, "MIA", "KC", 14, 21)
$games += ,@(398, 2023, 9, 45235, "TB", "HOU", 37, 39)
$games += ,@(399, 2023, 9, 45235, "LA", "GB", 3, 20)
$games += ,@(400, 2023, 9, 45235, "CHI", "NO", 17, 24)
$games += ,@(401, 2023, 9, 45235, "SEA", "BAL", 3, 37)
$games += ,@(402, 2023, 9, 45235, "MIN", "ATL", 31, 28)
$games += ,@(403, 2023, 9, 45235, "ARI", "CLE", 0, 27)
$games += ,@(404, 2023, 9, 45235, "WAS", "NE", 20, 17)
$games += ,@(405, 2023, 9, 45235, "IND", "CAR", 27, 13)
$games += ,@(406, 2023, 9, 45235, "NYG", "LV", 6, 30)
$games += ,@(407, 2023, 9, 45235, "DAL", "PHI", 23, 28)
$games += ,@(408, 2023, 9, 45235, "BUF", "CIN", 18, 24)
$games += ,@(409, 2023, 9, 45236, "LAC", "NYJ", 27, 6)
$games += ,@(410, 2023, 10, 45239, "CAR", "CHI", 13, 16)
$games += ,@(411, 2023, 10, 45242, "IND", "NE", 10, 6)
$games += ,@(412, 2023, 10, 45242, "NO", "MIN", 19, 27)
$games += ,@(413, 2023, 10, 45242, "HOU", "CIN", 30, 27)
$games += ,@(414, 2023, 10, 45242, "CLE", "BAL", 33, 31)
$games += ,@(415, 2023, 10, 45242, "GB", "PIT", 19, 23)
$games += ,@(416, 2023, 10, 45242, "TEN", "TB", 6, 20)
$games += ,@(417, 2023, 10, 45242, "SF", "JAX", 34, 3)
$games += ,@(418, 2023, 10, 45242, "DET", "LAC", 41, 38)
$games += ,@(419, 2023, 10, 45242, "ATL", "ARI", 23, 35)
$games += ,@(420, 2023, 10, 45242, "WAS", "SEA", 26, 29)
$games += ,@(421, 2023, 10, 45242, "NYG", "DAL", 17, 49)
$games += ,@(422, 2023, 10, 45242, "NYJ", "LV", 12, 16)
$games += ,@(423, 2023, 10, 45243, "DEN", "BUF", 24, 22)

foreach ($g in $games) {
    $r = $g[0]
    $ws.Cells.Item($r, 1).Value = $g[1]   # season
    $ws.Cells.Item($r, 2).Value = $g[2]   # week
    $ws.Cells.Item($r, 3).Value = $g[3]   # date
    $ws.Cells.Item($r, 4).Value = $g[4]   # team1 (away)
    $ws.Cells.Item($r, 5).Value = $g[5]   # team2 (home)
    $ws.Cells.Item($r, 6).Value = $g[6]   # score1
    $ws.Cells.Item($r, 7).Value = $g[7]   # score2
}

# ---------------------------------------------------------------
# 2. Column H ("home_team") is `=E<row>` for every game row. Rows
#    367-400 already carry that shared formula; rows 401-418 and
#    419-425 are new, so give each block its own formula fill.
# ---------------------------------------------------------------
$ws.Range("H401:H418").Formula = "=E401"
$ws.Range("H419:H425").Formula = "=E419"

# ---------------------------------------------------------------
# 3. Rows 424-440 get the season placeholder in column A only (424
#    and 425 already picked up the H formula above).
# ---------------------------------------------------------------
for ($r = 424; $r -le 440; $r++) {
    $ws.Cells.Item($r, 1).Value = 2023
}

# ---------------------------------------------------------------
# 4. Make the "games" sheet the active/selected tab again and
#    move the selection to where the last edit happened.
# ---------------------------------------------------------------
$ws.Activate()
$ws.Range("H423").Select()

